$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("InfoBasica")
$ws2 = $wb.Worksheets.Item("BuyBook")

# Update the email text displayed in A2 of both sheets (shared string used by both)
$ws1.Range("A2").Value = "pruebareto1001@yopmail.com"
$ws2.Range("A2").Value = "pruebareto1001@yopmail.com"

# Change active tab to BuyBook (sheet index 2, second sheet)
$ws2.Select()
